# Update "南宁-漫展信息.xlsx" with newly scraped convention-listing data.
#
# Affected sheets: "展览" (sheet1) and "全部类型" (sheet4).
#   - "展览": the oldest event (2024-08-24, id=88276) is dropped, every
#     remaining row shifts up by one, and one new event is appended at
#     the end (A1:I7 -> A1:I6).
#   - "全部类型": same shift/drop/append, but this sheet also keeps the
#     extra "莫西干人" live-music listing that "展览" does not carry
#     (A1:I8 -> A1:I7).
#
# NOTE: cells in column B hold plain date-like text (e.g. "2024-09-07").
# Excel's COM layer auto-recognizes such strings and would silently turn
# them into real date serials (with an automatically-applied date
# NumberFormat) if we just did `Cells.Item(r,2).Value = "..."`. To keep
# them as literal text - matching the original inlineStr content - we
# force the cell to Text format first, write the string, then restore
# the cell's style to "Normal" so no stray per-cell number format is
# left behind.

$wb = $excel.ActiveWorkbook

function Set-EventRow($ws, $r, $seqNo, $date, $title, $location, $timeRange, $wantCount, $minPrice, $link, $cover) {
    $ws.Cells.Item($r, 1).Value = $seqNo
    # Force column B to Text first so the date-shaped string ("2024-09-07")
    # is kept as literal text instead of being auto-converted to a real
    # date serial by Excel's input parser; restore the default style
    # afterwards so no stray per-cell number format lingers.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $date
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Value = $title
    $ws.Cells.Item($r, 4).Value = $location
    $ws.Cells.Item($r, 5).Value = $timeRange
    $ws.Cells.Item($r, 6).Value = $wantCount
    $ws.Cells.Item($r, 7).Value = $minPrice
    $ws.Cells.Item($r, 8).Value = $link
    $ws.Cells.Item($r, 9).Value = $cover
}

# ---------------------------------------------------------------------
# Sheet "展览": was A1:I7 (6 events), becomes A1:I6 (5 events).
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

Set-EventRow $wsExpo 2 1 "2024-09-07" "南宁·9.7国乙同人ONLY" `
    "南建路金砖茶城1号门 TZ· party" "2024.09.07 11:30-09.07 18:00" `
    36 109 `
    "https://show.bilibili.com/platform/detail.html?id=90932" `
    "//i2.hdslb.com/bfs/openplatform/202408/vSYiKkHQ1722860294516.jpeg"

Set-EventRow $wsExpo 3 2 "2024-09-15" "南宁·原神x星铁x绝区零同人ONLY3.0" `
    "亭洪路45号 百益上河城" "2024.09.15 10:00-09.15 17:00" `
    74 60 `
    "https://show.bilibili.com/platform/detail.html?id=90570" `
    "//i0.hdslb.com/bfs/openplatform/202408/sd7B5MV91723100089780.jpeg"

Set-EventRow $wsExpo 4 3 "2024-10-03" "南宁·2024良牙动漫秋季盛典（秋典）" `
    "民族大道106号 南宁国际会展中心" "2024.10.03 09:30-10.04 17:30" `
    1981 55 `
    "https://show.bilibili.com/platform/detail.html?id=90762" `
    "//i0.hdslb.com/bfs/openplatform/202408/njVhnU591723691579900.jpeg"

Set-EventRow $wsExpo 5 4 "2024-10-03" "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini" `
    "南宁国际会展中心  南宁国际会展中心" "2024.10.03 09:30-10.04 17:30" `
    149 55 `
    "https://show.bilibili.com/platform/detail.html?id=91043" `
    "//i2.hdslb.com/bfs/openplatform/202408/jEAI96Ev1724123680899.jpeg"

Set-EventRow $wsExpo 6 5 "2024-11-02" "南宁·万圣漫控嘉年华10" `
    "亭洪路45号 百益上河城" "2024.11.02 11:00-11.03 22:00" `
    349 50 `
    "https://show.bilibili.com/platform/detail.html?id=87820" `
    "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"

# Old row 7 (A7:I7) no longer exists in the new data; remove it so the
# used range shrinks from A1:I7 to A1:I6, matching the new <dimension>.
$wsExpo.Range("A7:I7").Delete()

# ---------------------------------------------------------------------
# Sheet "全部类型": was A1:I8 (7 events), becomes A1:I7 (6 events).
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

Set-EventRow $wsAll 2 1 "2024-09-07" "南宁·9.7国乙同人ONLY" `
    "南建路金砖茶城1号门 TZ· party" "2024.09.07 11:30-09.07 18:00" `
    36 109 `
    "https://show.bilibili.com/platform/detail.html?id=90932" `
    "//i2.hdslb.com/bfs/openplatform/202408/vSYiKkHQ1722860294516.jpeg"

Set-EventRow $wsAll 3 2 "2024-09-15" "南宁·原神x星铁x绝区零同人ONLY3.0" `
    "亭洪路45号 百益上河城" "2024.09.15 10:00-09.15 17:00" `
    74 60 `
    "https://show.bilibili.com/platform/detail.html?id=90570" `
    "//i0.hdslb.com/bfs/openplatform/202408/sd7B5MV91723100089780.jpeg"

Set-EventRow $wsAll 4 3 "2024-10-03" "南宁·2024良牙动漫秋季盛典（秋典）" `
    "民族大道106号 南宁国际会展中心" "2024.10.03 09:30-10.04 17:30" `
    1981 55 `
    "https://show.bilibili.com/platform/detail.html?id=90762" `
    "//i0.hdslb.com/bfs/openplatform/202408/njVhnU591723691579900.jpeg"

Set-EventRow $wsAll 5 4 "2024-10-03" "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini" `
    "南宁国际会展中心  南宁国际会展中心" "2024.10.03 09:30-10.04 17:30" `
    149 55 `
    "https://show.bilibili.com/platform/detail.html?id=91043" `
    "//i2.hdslb.com/bfs/openplatform/202408/jEAI96Ev1724123680899.jpeg"

Set-EventRow $wsAll 6 5 "2024-10-04" "南宁·《最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会》" `
    "福建园街道星光大道4号 南宁剧场" "2024.10.04 20:00-10.04 21:30" `
    6 100 `
    "https://show.bilibili.com/platform/detail.html?id=89039" `
    "//i0.hdslb.com/bfs/openplatform/202407/dudapgjU1720595605665.jpeg"

Set-EventRow $wsAll 7 6 "2024-11-02" "南宁·万圣漫控嘉年华10" `
    "亭洪路45号 百益上河城" "2024.11.02 11:00-11.03 22:00" `
    349 50 `
    "https://show.bilibili.com/platform/detail.html?id=87820" `
    "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"

# Old row 8 (A8:I8) no longer exists in the new data; remove it so the
# used range shrinks from A1:I8 to A1:I7, matching the new <dimension>.
$wsAll.Range("A8:I8").Delete()
